# Auto-generated edit script applying the cryptos.xlsx diff
# (GitHub Actions crypto price refresh, Fri Mar 29 06:41:29 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.396.17'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '3.565.47'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''607.68'
$ws.Range("E5").Value = '  +3.53%  '
$ws.Range("D6").Value = '''186.39'
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").Value = '3.558.72'
$ws.Range("E7").Value = '  +0.39%  '
$ws.Range("D8").Value = '''0.619'
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").Value = '''0.217'
$ws.Range("E10").Value = '  +9.15%  '
$ws.Range("D11").Value = '''0.646'
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("D12").Value = '''53.95'
$ws.Range("E12").Value = '  -0.80%  '
$ws.Range("D13").Value = '''0.0000311'
$ws.Range("E13").Value = '  +1.65%  '
$ws.Range("D14").Value = '''9.56'
$ws.Range("E14").Value = '  +0.69%  '
$ws.Range("D15").Value = '4.129.27'
$ws.Range("E15").Value = '  +0.15%  '
$ws.Range("D16").Value = '70.461.29'
$ws.Range("E16").Value = '  +0.57%  '
# Row 17: coin re-ranked
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").Value = '''12.78'
$ws.Range("E17").Value = '  +2.54%  '
# Row 18: coin re-ranked
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '''19.06'
$ws.Range("E18").Value = '  -1.70%  '
# Row 19: coin re-ranked
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '''582.37'
$ws.Range("E19").Value = '  +8.27%  '
# Row 20: coin re-ranked
$ws.Range("B20").Value = 'WrappedEther'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D20").Value = '3.539.46'
$ws.Range("E20").Value = '  -0.29%  '
$ws.Range("E21").Value = '  +0.42%  '
$ws.Range("D22").Value = '''0.996'
$ws.Range("E22").Value = '  -1.88%  '
$ws.Range("D23").Value = '''17.37'
$ws.Range("E23").Value = '  -3.97%  '
$ws.Range("D24").Value = '''4.73'
$ws.Range("E24").Value = '  +2.05%  '
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("D26").Value = '''94.37'
$ws.Range("E26").Value = '  -1.00%  '
$ws.Range("D27").Value = '''2.95'
$ws.Range("E27").Value = '  -1.33%  '
$ws.Range("D28").Value = '''10.94'
$ws.Range("E28").Value = '  -2.59%  '
$ws.Range("D29").Value = '''9.44'
$ws.Range("E29").Value = '  +3.39%  '
$ws.Range("D30").Value = '''32.36'
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("D31").Value = '''7.09'
$ws.Range("E31").Value = '  -3.27%  '
$ws.Range("D32").Value = '''12.26'
$ws.Range("D33").Value = '''0.115'
$ws.Range("E33").Value = '  +1.15%  '
$ws.Range("D34").Value = '''63.73'
$ws.Range("E34").Value = '  -1.98%  '
$ws.Range("D35").Value = '''3.71'
$ws.Range("E35").Value = '  +20.26%  '
$ws.Range("D36").Value = '''3.21'
$ws.Range("E36").Value = '  -0.39%  '
$ws.Range("D37").Value = '''531.24'
$ws.Range("E37").Value = '  -3.54%  '
$ws.Range("D38").Value = '''0.405'
$ws.Range("E38").Value = '  -2.24%  '
$ws.Range("D39").Value = '''1.00'
$ws.Range("E39").Value = '  +0.10%  '
$ws.Range("D40").Value = '''37.50'
$ws.Range("E40").Value = '  -2.64%  '
$ws.Range("D41").Value = '0.0₃0789'
$ws.Range("E41").Value = '  +2.90%  '
$ws.Range("D42").Value = '3.528.25'
$ws.Range("E42").Value = '  +5.37%  '
$ws.Range("E43").Value = '  +4.43%  '
$ws.Range("E44").Value = '  +1.41%  '
$ws.Range("D45").Value = '''0.0461'
$ws.Range("E45").Value = '  +4.07%  '
$ws.Range("D46").Value = '''3.48'
$ws.Range("E46").Value = '  -3.62%  '
$ws.Range("E47").Value = '  -1.81%  '
$ws.Range("D49").Value = '''9.22'
$ws.Range("E49").Value = '  +0.28%  '
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("D51").Value = '''135.53'
$ws.Range("E51").Value = '  -1.31%  '
